$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'48.322.00"
$ws.Range("E2").Value = "  +1.43%  "

$ws.Range("D3").Value = "'2.508.74"

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'321.79"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").Value = "'108.32"
$ws.Range("E6").Value = "  -0.92%  "

$ws.Range("E7").Value = "  +1.16%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -0.51%  "

$ws.Range("D10").Value = "'39.85"
$ws.Range("E10").Value = "  +0.71%  "

$ws.Range("E11").Value = "  +8.98%  "

$ws.Range("E12").Value = "  +0.95%  "

$ws.Range("E13").Value = "  -0.11%  "

$ws.Range("D14").Value = "'7.20"
$ws.Range("E14").Value = "  -0.16%  "

$ws.Range("D15").Value = "'2.901.44"
$ws.Range("E15").Value = "  +0.64%  "

$ws.Range("D16").Value = "'2.508.39"
$ws.Range("E16").Value = "  +0.55%  "

$ws.Range("E17").Value = "  -0.46%  "

$ws.Range("D18").Value = "'48.155.77"
$ws.Range("E18").Value = "  +1.37%  "

$ws.Range("D19").Value = "'13.11"
$ws.Range("E19").Value = "  -2.29%  "

$ws.Range("D20").Value = "'6.81"
$ws.Range("E20").Value = "  +2.43%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.0₃0946"
$ws.Range("E21").Value = "  +0.35%  "

$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").Value = "'2.79"
$ws.Range("E22").Value = "  +1.08%  "

$ws.Range("D23").Value = "'281.26"
$ws.Range("E23").Value = "  +13.83%  "

$ws.Range("D24").Value = "'72.44"
$ws.Range("E24").Value = "  +2.45%  "

$ws.Range("D25").Value = "'2.56"
$ws.Range("E25").Value = "  +0.28%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").Value = "'25.78"
$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("E28").Value = "  -1.11%  "

$ws.Range("D29").Value = "'9.80"
$ws.Range("E29").Value = "  -2.04%  "

$ws.Range("E30").Value = "  +0.89%  "

$ws.Range("D31").Value = "'35.40"
$ws.Range("E31").Value = "  +1.60%  "

$ws.Range("D32").Value = "'49.44"
$ws.Range("E32").Value = "  -0.98%  "

$ws.Range("D33").Value = "'19.68"
$ws.Range("E33").Value = "  -3.59%  "

$ws.Range("E34").Value = "  +0.52%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("E36").Value = "  -0.64%  "

$ws.Range("E37").Value = "  -0.45%  "

$ws.Range("E38").Value = "  -1.38%  "

$ws.Range("E39").Value = "  -0.96%  "

$ws.Range("E40").Value = "  -0.14%  "

$ws.Range("D41").Value = "'121.70"
$ws.Range("E41").Value = "  +2.30%  "

$ws.Range("E42").Value = "  -0.17%  "

$ws.Range("D43").Value = "'21.46"
$ws.Range("E43").Value = "  -4.50%  "

$ws.Range("E44").Value = "  +1.80%  "

$ws.Range("D45").Value = "'2.020.02"
$ws.Range("E45").Value = "  +1.07%  "

$ws.Range("E46").Value = "  +4.31%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'1.85"
$ws.Range("E47").Value = "  +2.82%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'2.00"
$ws.Range("E48").Value = "  -2.59%  "

$ws.Range("E49").Value = "  -0.53%  "

$ws.Range("D50").Value = "'5.19"
$ws.Range("E50").Value = "  -0.86%  "

$ws.Range("D51").Value = "'80.82"
$ws.Range("E51").Value = "  +4.01%  "
